$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4002.2
$ws.Range("I113").Value = 2002.5
$ws.Range("K113").Value = 2002.5
$ws.Range("M113").Value = 1251.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2244.9512
$ws.Range("I132").Value = 1969.3846
$ws.Range("J132").Value = 2722.6
$ws.Range("K132").Value = 5908.1538
$ws.Range("L132").Value = 8167.799999999999
$ws.Range("M132").Value = -3378.1538
$ws.Range("N132").Value = -13227.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4141.7646
$ws.Range("I141").Value = 3887.1428
$ws.Range("J141").Value = 5330
$ws.Range("K141").Value = 11661.4284
$ws.Range("L141").Value = 15990
$ws.Range("M141").Value = -6481.428400000001
$ws.Range("N141").Value = -26350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5733.61
$ws.Range("I32").Value = 4359.222
$ws.Range("K32").Value = 4359.222
$ws.Range("M32").Value = -4072.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3033116
$ws.Range("I45").Value = 6062866
$ws.Range("J45").Value = 3365.8
$ws.Range("K45").Value = 6062866
$ws.Range("L45").Value = 3365.8
$ws.Range("M45").Value = -6062489
$ws.Range("N45").Value = -4119.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6166.467
$ws.Range("I74").Value = 3380.476
$ws.Range("J74").Value = 12667.111
$ws.Range("K74").Value = 3380.476
$ws.Range("L74").Value = 12667.111
$ws.Range("M74").Value = -2506.476
$ws.Range("N74").Value = -14415.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6166.467
$ws.Range("I77").Value = 3380.476
$ws.Range("J77").Value = 12667.111
$ws.Range("K77").Value = 16902.38
$ws.Range("L77").Value = 63335.55500000001
$ws.Range("M77").Value = -12534.38
$ws.Range("N77").Value = -72071.55500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3170.5
$ws.Range("I122").Value = 4919
$ws.Range("J122").Value = 2121.4
$ws.Range("K122").Value = 14757
$ws.Range("L122").Value = 6364.200000000001
$ws.Range("M122").Value = -12307
$ws.Range("N122").Value = -11264.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 41780
$ws.Range("J52").Value = 41780
$ws.Range("L52").Value = 41780
$ws.Range("N52").Value = -42306

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 57000
$ws.Range("J57").Value = 57000
$ws.Range("L57").Value = 57000
$ws.Range("N57").Value = -58440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H121").Value = 41780
$ws.Range("J121").Value = 41780
$ws.Range("L121").Value = 41780
$ws.Range("N121").Value = -45274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3224.0981
$ws.Range("I134").Value = 3088.58
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 9265.74
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -6730.74
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 57000
$ws.Range("J136").Value = 57000
$ws.Range("L136").Value = 57000
$ws.Range("N136").Value = -67200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1494.3846
$ws.Range("I94").Value = 1485.5
$ws.Range("J94").Value = 1502
$ws.Range("K94").Value = 1485.5
$ws.Range("L94").Value = 1502
$ws.Range("M94").Value = -1034.5
$ws.Range("N94").Value = -2404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 69000
$ws.Range("J97").Value = 69000
$ws.Range("L97").Value = 69000
$ws.Range("N97").Value = -70982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 11202.77
$ws.Range("I122").Value = 3550.1667
$ws.Range("K122").Value = 10650.5001
$ws.Range("M122").Value = -8200.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10405.143
$ws.Range("I132").Value = 12547.4
$ws.Range("J132").Value = 5049.5
$ws.Range("K132").Value = 37642.2
$ws.Range("L132").Value = 15148.5
$ws.Range("M132").Value = -35112.2
$ws.Range("N132").Value = -20208.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 749.29
$ws.Range("I113").Value = 753.2941
$ws.Range("J113").Value = 726.6
$ws.Range("K113").Value = 2259.8823
$ws.Range("L113").Value = 2179.8
$ws.Range("M113").Value = -89.88229999999976
$ws.Range("N113").Value = -6519.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2083
$ws.Range("I132").Value = 2737.8
$ws.Range("J132").Value = 1673.75
$ws.Range("K132").Value = 24640.2
$ws.Range("L132").Value = 15063.75
$ws.Range("M132").Value = -22110.2
$ws.Range("N132").Value = -20123.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4773.4
$ws.Range("I113").Value = 5877.75
$ws.Range("J113").Value = 4037.1667
$ws.Range("K113").Value = 5877.75
$ws.Range("L113").Value = 4037.1667
$ws.Range("M113").Value = -3707.75
$ws.Range("N113").Value = -8377.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 18200
$ws.Range("I122").Value = 26050
$ws.Range("K122").Value = 78150
$ws.Range("M122").Value = -75700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 29179.5
$ws.Range("J123").Value = 29179.5
$ws.Range("L123").Value = 29179.5
$ws.Range("N123").Value = -34079.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 17914.79
$ws.Range("I61").Value = 25013.54
$ws.Range("J61").Value = 2534.1667
$ws.Range("K61").Value = 25013.54
$ws.Range("L61").Value = 2534.1667
$ws.Range("M61").Value = -24811.54
$ws.Range("N61").Value = -2938.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1433.3334
$ws.Range("I68").Value = 1457.1428
$ws.Range("J68").Value = 1350
$ws.Range("K68").Value = 1457.1428
$ws.Range("L68").Value = 1350
$ws.Range("M68").Value = -708.1428000000001
$ws.Range("N68").Value = -2848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1433.3334
$ws.Range("I71").Value = 1457.1428
$ws.Range("J71").Value = 1350
$ws.Range("K71").Value = 7285.714
$ws.Range("L71").Value = 6750
$ws.Range("M71").Value = -3541.714
$ws.Range("N71").Value = -14238

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 17914.79
$ws.Range("I113").Value = 25013.54
$ws.Range("J113").Value = 2534.1667
$ws.Range("K113").Value = 25013.54
$ws.Range("L113").Value = 2534.1667
$ws.Range("M113").Value = -22843.54
$ws.Range("N113").Value = -6874.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2908.5833
$ws.Range("J4").Value = 2991.182
$ws.Range("L4").Value = 2991.182
$ws.Range("N4").Value = -3217.182
